$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the affected rows
$ws.Range("F2").Value = -5
$ws.Range("F3").Value = -4
$ws.Range("F5").Value = 8
$ws.Range("F8").Value = -3
$ws.Range("F10").Value = -3
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = -6
$ws.Range("F16").Value = 7
$ws.Range("F17").Value = -8
$ws.Range("F18").Value = 2
$ws.Range("F20").Value = 0
$ws.Range("F21").Value = -5
$ws.Range("F22").Value = -5
